$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) for rows 2-22 from 45207 to 45208 (one day later)
for ($row = 2; $row -le 22; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45207) {
        $cell.Value2 = 45208
    }
}
